# corrección de errores cuando no hay cliente asociado
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 61: "Error en iva cuando consumidor final en ticket" -> mark as 100% done (C61 = 1)
$ws.Range("C61").Value = 1
$ws.Range("C61").NumberFormat = "0%"

# Row 68: "Cuando no hay cliente asociado no setea en ventaDTO el campo tipoTicket" -> mark as 100% done (C68 = 1)
$ws.Range("C68").Value = 1
$ws.Range("C68").NumberFormat = "0%"

# Row 69: "Setear corte z en el ticket" -> status "en proceso" (C69)
$ws.Range("C69").Value = "en proceso"

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("C70").Select()
